$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.690.30"
$ws.Range("E2").Value = "  +0.91%  "
$ws.Range("D3").Value = "2.244.87"
$ws.Range("E3").Value = "  +0.12%  "
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "306.38"
$ws.Range("E5").Value = "  -0.33%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "94.25"
$ws.Range("E6").Value = "  -0.28%  "
$ws.Range("E7").Value = "  -0.29%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.514"
$ws.Range("E9").Value = "  -2.15%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.74"
$ws.Range("E10").Value = "  -0.16%  "
$ws.Range("E11").Value = "  -1.50%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.17"
$ws.Range("E12").Value = "  -0.35%  "
$ws.Range("E13").Value = "  -0.01%  "
$ws.Range("D14").Value = "2.586.91"
$ws.Range("E14").Value = "  +0.19%  "
$ws.Range("D15").Value = "2.344.67"
$ws.Range("E15").Value = "  +4.63%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.829"
$ws.Range("E16").Value = "  -0.25%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.52"
$ws.Range("E17").Value = "  +0.03%  "
$ws.Range("D18").Value = "44.446.01"
$ws.Range("E18").Value = "  +0.95%  "
$ws.Range("D19").Value = "0.0₃0933"
$ws.Range("E19").Value = "  -3.18%  "
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.17"
$ws.Range("E20").Value = "  -3.65%  "
$ws.Range("B21").Value = "InternetComputer(DFINITY)"
$ws.Range("C21").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.73"
$ws.Range("E21").Value = "  -3.40%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "65.31"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "237.58"
$ws.Range("E23").Value = "  -0.31%  "
$ws.Range("E24").Value = "  -0.25%  "
$ws.Range("E25").Value = "  -1.56%  "
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.29"
$ws.Range("E27").Value = "  +3.62%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.76"
$ws.Range("E28").Value = "  -1.66%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "36.89"
$ws.Range("E29").Value = "  -4.28%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.93"
$ws.Range("E30").Value = "  -0.48%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.83"
$ws.Range("E31").Value = "  -0.41%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "148.61"
$ws.Range("E32").Value = "  -2.91%  "
$ws.Range("B33").Value = "WEMIXToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.62"
$ws.Range("E33").Value = "  +0.34%  "
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0779"
$ws.Range("E34").Value = "  -2.04%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.17"
$ws.Range("E35").Value = "  +1.00%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.108"
$ws.Range("E36").Value = "  +1.44%  "
$ws.Range("E37").Value = "  -2.08%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.85"
$ws.Range("E38").Value = "  +4.86%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "15.04"
$ws.Range("E39").Value = "  +4.97%  "
$ws.Range("E40").Value = "  -5.03%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.76"
$ws.Range("E41").Value = "  -1.37%  "
$ws.Range("E42").Value = "  -0.68%  "
$ws.Range("E43").Value = "  -0.01%  "
$ws.Range("D44").Value = "1.809.49"
$ws.Range("E44").Value = "  +3.29%  "
$ws.Range("E45").Value = "  +13.54%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "81.68"
$ws.Range("E46").Value = "  -1.58%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.186"
$ws.Range("E47").Value = "  -2.58%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "97.78"
$ws.Range("E48").Value = "  -2.10%  "
$ws.Range("B49").Value = "THORChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.82"
$ws.Range("E49").Value = "  -2.27%  "
$ws.Range("B50").Value = "ordi"
$ws.Range("C50").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "68.87"
$ws.Range("E50").Value = "  +2.17%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "53.80"
$ws.Range("E51").Value = "  -1.75%  "
